# Weekly update: insert a new week's worth of price data (rows 338-340)
# right before the existing 2021-06-10 block, shifting everything else down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows at row 338 (existing rows 338+ shift down by 3).
$ws.Range("A338:R340").EntireRow.Insert()

# Row 338: Tomate, Larga vida, Primera
$ws.Cells.Item(338, 1).Value = 2
$ws.Cells.Item(338, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(338, 3).Value = "Coquimbo"
$ws.Cells.Item(338, 4).Value = 44461
$ws.Cells.Item(338, 5).Value = 4
$ws.Cells.Item(338, 6).Value = 100112020
$ws.Cells.Item(338, 7).Value = "Tomate"
$ws.Cells.Item(338, 8).Value = "Larga vida"
$ws.Cells.Item(338, 9).Value = "Primera"
$ws.Cells.Item(338, 10).Value = 2400
$ws.Cells.Item(338, 11).Value = 11000
$ws.Cells.Item(338, 12).Value = 12000
$ws.Cells.Item(338, 13).Value = 11500
$ws.Cells.Item(338, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(338, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(338, 16).Value = 639
$ws.Cells.Item(338, 17).Value = 18
$ws.Cells.Item(338, 18).Value = "Hortaliza"

# Row 339: Tomate, Larga vida, Segunda
$ws.Cells.Item(339, 1).Value = 2
$ws.Cells.Item(339, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(339, 3).Value = "Coquimbo"
$ws.Cells.Item(339, 4).Value = 44461
$ws.Cells.Item(339, 5).Value = 4
$ws.Cells.Item(339, 6).Value = 100112020
$ws.Cells.Item(339, 7).Value = "Tomate"
$ws.Cells.Item(339, 8).Value = "Larga vida"
$ws.Cells.Item(339, 9).Value = "Segunda"
$ws.Cells.Item(339, 10).Value = 2000
$ws.Cells.Item(339, 11).Value = 9000
$ws.Cells.Item(339, 12).Value = 10000
$ws.Cells.Item(339, 13).Value = 9500
$ws.Cells.Item(339, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(339, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(339, 16).Value = 528
$ws.Cells.Item(339, 17).Value = 18
$ws.Cells.Item(339, 18).Value = "Hortaliza"

# Row 340: Tomate, Larga vida, Tercera
$ws.Cells.Item(340, 1).Value = 2
$ws.Cells.Item(340, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(340, 3).Value = "Coquimbo"
$ws.Cells.Item(340, 4).Value = 44461
$ws.Cells.Item(340, 5).Value = 4
$ws.Cells.Item(340, 6).Value = 100112020
$ws.Cells.Item(340, 7).Value = "Tomate"
$ws.Cells.Item(340, 8).Value = "Larga vida"
$ws.Cells.Item(340, 9).Value = "Tercera"
$ws.Cells.Item(340, 10).Value = 1600
$ws.Cells.Item(340, 11).Value = 7000
$ws.Cells.Item(340, 12).Value = 8000
$ws.Cells.Item(340, 13).Value = 7500
$ws.Cells.Item(340, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(340, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(340, 16).Value = 417
$ws.Cells.Item(340, 17).Value = 18
$ws.Cells.Item(340, 18).Value = "Hortaliza"
